$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07796894984218658
$ws.Range("D2").Value = 0.1911874935925047
$ws.Range("G2").Value = 0.12469127785007
$ws.Range("H2").Value = 0.99
